# Applies updated vm_pu values (case with 380 kV) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B,C,D,E,F,I,J,K,L,M map to spreadsheet column indices 2,3,4,5,6,9,10,11,12,13
$cols = @(2,3,4,5,6,9,10,11,12,13)

# New values for rows 2..25 (one row per array entry), in the column order above
$data = @(
    (1.02,1.040667427109522,1.042112379478266,1.038958755264217,1.039424397813519,1.041318843429399,1.045752794676575,1.044889987376645,1.041745306298796,1.04220962465886),
    (1.02,1.041984716831792,1.043114880628447,1.040095706977138,1.041371098644692,1.04177066314492,1.046713961725308,1.045702751054996,1.042691515913576,1.043963548044885),
    (1.02,1.042835418703603,1.043762072289577,1.040830049940092,1.042628803438122,1.042060812318833,1.047333774353137,1.046226564695277,1.043301860991411,1.045096119351156),
    (1.02,1.043192660025108,1.044033798447362,1.041138451802652,1.043157092659726,1.042182265077159,1.047593839600363,1.04644627771751,1.043557996199931,1.045571706786447),
    (1.02,1.043252619519143,1.044079401878086,1.041190215453555,1.043245768790876,1.042202626766662,1.047637476295208,1.046483139406332,1.043600975952237,1.045651528403996),
    (1.02,1.042840193718516,1.043765704490335,1.040834172055963,1.042635864214066,1.042062437239813,1.047337251331653,1.046229502463949,1.043305285260077,1.045102476300952),
    (1.02,1.04111296171713,1.042451491333662,1.039343274111386,1.040082705246498,1.041471996936039,1.046078069050887,1.045165101793126,1.042065481842463,1.042802862676393),
    (1.02,1.038056246265958,1.040124043781676,1.036705623539756,1.035568184804161,1.04041452839692,1.043842708239259,1.043273226828707,1.039865901199612,1.038732171862284),
    (1.02,1.036009196029881,1.038564310620093,1.034939808845736,1.032547133898603,1.039697919846801,1.042341032273725,1.04200077216447,1.038389199639641,1.036005114688603),
    (1.02,1.035120513878519,1.037886944626094,1.03417336824494,1.031236080923096,1.039384822669034,1.041688003951532,1.041447067030392,1.037747256663492,1.034820935432875),
    (1.02,1.034790064527948,1.037635036067419,1.033888396607717,1.030748641132394,1.039268100156659,1.041445014134921,1.041240981744084,1.037508425649845,1.03438055968523),
    (1.02,1.034860963153893,1.037689085184315,1.033949536831208,1.030853219562464,1.039293156787416,1.041497155657988,1.041285206582548,1.037559673223234,1.034475045518973),
    (1.02,1.035093206086204,1.037866128033265,1.034149818179457,1.031195798411578,1.039375183025125,1.041667927065723,1.041430040443885,1.037727522702258,1.034784544514906),
    (1.02,1.035236251565663,1.037975169438346,1.03427318063917,1.031406811537664,1.039425665771797,1.041773088320244,1.041519222274706,1.037830889134815,1.034975167820904),
    (1.02,1.036068125859031,1.038609222748328,1.034990635862002,1.032634081160089,1.039718639767267,1.042384312224516,1.042037461909385,1.038431749646674,1.036083632898729),
    (1.02,1.036589318114724,1.039006410655935,1.035440181838919,1.03340312112418,1.039901662282845,1.042766964608901,1.04236180694764,1.03880797429517,1.036778036665955),
    (1.02,1.036893099784541,1.039237891892187,1.035702218151954,1.033851408800259,1.040008146205442,1.042989890253042,1.042550729469562,1.039027176891571,1.037182748945907),
    (1.02,1.036996644255387,1.039316788653007,1.035791536059202,1.03400421630555,1.040044408776438,1.043065856700581,1.042615102810815,1.039101878296382,1.037320691299053),
    (1.02,1.036533422034409,1.038963816027171,1.035391968087813,1.033320639542132,1.039882053653762,1.042725937503237,1.042327034980832,1.038767634117387,1.036703567107399),
    (1.02,1.035024826149077,1.03781400178179,1.034090848135658,1.031094930221273,1.039351040082806,1.041617650954882,1.041387401968056,1.037678105885246,1.034693419228411),
    (1.02,1.034074264189099,1.037089302441055,1.033271150675078,1.029692889174026,1.039014714278976,1.040918359145933,1.040794216293978,1.036990846843737,1.033426549442883),
    (1.02,1.034578371972071,1.037473648474435,1.033705844613078,1.030436394815041,1.039193240965555,1.041289303077602,1.04110890462836,1.037355389200487,1.034098431858138),
    (1.02,1.036558679749094,1.038983063298488,1.035413754350675,1.033357910282996,1.039890914783202,1.042744476716958,1.042342747746373,1.038785862857248,1.036737217687711),
    (1.02,1.038848077759144,1.040727152142898,1.03738879653661,1.036737232004906,1.040689946237862,1.044422593942828,1.043764276874556,1.040436341041456,1.039786818744906)
)

$rowStart = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $rowStart + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $vals[$j]
    }
}
